$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.554.58'
$ws.Range("E2").Value = '  +1.10%  '
$ws.Range("D3").Value = '1.791.08'
$ws.Range("E3").Value = '  -0.61%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.006'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.22%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '329.12'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.74%  '
$ws.Range("E6").Value = '  +0.25%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4390'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.47%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3737'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +5.76%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '45.45'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.21%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07589'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.38%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.131'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.32%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '22.59'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.40%  '
$ws.Range("E13").Value = '  +0.19%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.211'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.36%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.479'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.27%  '
$ws.Range("D16").Value = '1.793.69'
$ws.Range("E16").Value = '  -0.27%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001087'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.18%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06696'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.25%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '80.39'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.33%  '
$ws.Range("E20").Value = '  +0.46%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.49'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.69%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.201'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.04%  '
$ws.Range("D23").Value = '28.552.85'
$ws.Range("E23").Value = '  +1.15%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.68'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.27%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.436'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.56%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '20.36'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.49%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '152.74'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.60%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.329'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.47%  '
$ws.Range("D29").Value = '1.999.27'
$ws.Range("E29").Value = '  -0.30%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.300'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.01%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '130.46'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.77%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.977'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.34%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.780'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.40%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.09238'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.30%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.2240'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.62%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '12.07'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.77%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06250'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.28%  '
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02315'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.58%  '
$ws.Range("B39").Value = 'InternetComputer(DFINITY)'
$ws.Range("C39").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.188'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.05%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.6569'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.22%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.197'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.40%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.423'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.07%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '7.985'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.10%  '
$ws.Range("E44").Value = '  +0.30%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.88'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.34%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6065'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.61%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.811'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.45%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '127.34'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.40%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.007'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.07%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06998'
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.136'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.41%  '
